# Update Betfair Back/Lay odds for the games listed on 2025-10-14.
# Values below reflect the latest odds refresh for rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Paysandu vs Remo)
$ws.Range("F2").Value = 3.05
$ws.Range("G2").Value = 3.45
$ws.Range("I2").Value = 2.72
$ws.Range("J2").Value = 3.1
$ws.Range("Q2").Value = 2.26
$ws.Range("V2").Value = 1.58

# Row 3 (Chapecoense vs Botafogo SP)
$ws.Range("F3").Value = 1.63
$ws.Range("G3").Value = 1.72
$ws.Range("J3").Value = 3.75
$ws.Range("L3").Value = 1.47
$ws.Range("N3").Value = 3.1
$ws.Range("P3").Value = 1.73
$ws.Range("R3").Value = 1.27
$ws.Range("U3").Value = 1.77
$ws.Range("W3").Value = 2.38
$ws.Range("AC3").Value = 10.5
$ws.Range("AG3").Value = 12
$ws.Range("AJ3").Value = 19.5
$ws.Range("AN3").Value = 15.5
